# Updated cryptos list on Sun Sep 10 18:40:21 UTC 2023 with GitHub Actions
# Applies per-cell Price / Volume(1h) updates (and the ShibaInu/Litecoin row swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: Price (column D) values are forced to text with a leading single-quote
# so Excel does not auto-convert numeric-looking strings (e.g. "213.76") into
# real numbers, matching the inline-string cell type used in the source file.

$ws.Range("D2").Value = '''25.932.45'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = '''1.624.65'
$ws.Range("E3").Value = '  -0.99%  '
$ws.Range("D5").Value = '''213.76'
$ws.Range("E5").Value = '  -0.95%  '
$ws.Range("D6").Value = '''0.502'
$ws.Range("E6").Value = '  -0.89%  '
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("D8").Value = '''0.249'
$ws.Range("E8").Value = '  -2.28%  '
$ws.Range("D9").Value = '''0.0616'
$ws.Range("E9").Value = '  -3.43%  '
$ws.Range("D10").Value = '''18.22'
$ws.Range("E10").Value = '  -6.69%  '
$ws.Range("E11").Value = '  -0.94%  '
$ws.Range("D12").Value = '''1.849.56'
$ws.Range("E12").Value = '  -1.03%  '
$ws.Range("D13").Value = '''1.637.65'
$ws.Range("E13").Value = '  -0.33%  '
$ws.Range("E14").Value = '  -1.84%  '
$ws.Range("E15").Value = '  -3.63%  '
$ws.Range("D16").Value = '''25.926.20'
$ws.Range("E16").Value = '  -0.61%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = '''61.19'
$ws.Range("E17").Value = '  -3.39%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '''0.0₃0734'
$ws.Range("E18").Value = '  -3.59%  '
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("D20").Value = '''191.66'
$ws.Range("E20").Value = '  -1.38%  '
$ws.Range("D21").Value = '''4.23'
$ws.Range("E21").Value = '  -2.89%  '
$ws.Range("D22").Value = '''9.56'
$ws.Range("E22").Value = '  -3.59%  '
$ws.Range("D23").Value = '''6.06'
$ws.Range("E23").Value = '  -2.20%  '
$ws.Range("D24").Value = '''0.133'
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("D25").Value = '''143.47'
$ws.Range("E25").Value = '  +0.35%  '
$ws.Range("E26").Value = '  +0.31%  '
$ws.Range("E27").Value = '  -2.33%  '
$ws.Range("D28").Value = '''6.70'
$ws.Range("E28").Value = '  -2.39%  '
$ws.Range("D29").Value = '''15.16'
$ws.Range("E29").Value = '  -2.24%  '
$ws.Range("E30").Value = '  -1.51%  '
$ws.Range("E31").Value = '  -2.70%  '
$ws.Range("D32").Value = '''3.12'
$ws.Range("E32").Value = '  -3.89%  '
$ws.Range("D33").Value = '''3.11'
$ws.Range("E33").Value = '  -5.36%  '
$ws.Range("E34").Value = '  -2.88%  '
$ws.Range("E35").Value = '  -2.40%  '
$ws.Range("D36").Value = '''1.118.89'
$ws.Range("E36").Value = '  -0.87%  '
$ws.Range("E37").Value = '  -6.42%  '
$ws.Range("D38").Value = '''2.43'
$ws.Range("E38").Value = '  -1.41%  '
$ws.Range("E39").Value = '  -3.61%  '
$ws.Range("E40").Value = '  -2.04%  '
$ws.Range("E41").Value = '  -1.16%  '
$ws.Range("E42").Value = '  -3.68%  '
$ws.Range("D43").Value = '''1.760.60'
$ws.Range("E43").Value = '  -0.96%  '
$ws.Range("E44").Value = '  -5.61%  '
$ws.Range("D45").Value = '''0.0₆0115'
$ws.Range("E45").Value = '  -2.32%  '
$ws.Range("E46").Value = '  +1.16%  '
$ws.Range("D47").Value = '''54.46'
$ws.Range("E47").Value = '  -3.34%  '
$ws.Range("D48").Value = '''1.46'
$ws.Range("E48").Value = '  -0.92%  '
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("D51").Value = '''7.47'
$ws.Range("E51").Value = '  -3.78%  '
